$p = $ppt.ActivePresentation
$s = $p.Slides.Item(3)
$shp = $s.Shapes.Item(18)
$tr = $shp.TextFrame.TextRange

$tr.Text = " etc."
$tr.Font.Size = 12
$tr.Font.Bold = $true
$tr.LanguageID = "en-US"

$tr.InsertBefore(":.,:-")
$tr2 = $shp.TextFrame.TextRange
$tr2.LanguageID = "mr-IN"

$tr3 = $shp.TextFrame.TextRange
$tr3.InsertBefore("cut adjectives,  trim ")
$tr4 = $shp.TextFrame.TextRange
$tr4.LanguageID = "en-US"

# Try whole-range Font BEFORE doing the paragraph split
$trFull = $shp.TextFrame.TextRange
$trFull.Font.Size = 12
$trFull.Font.Bold = $true
Write-Output "pre-split: [$($trFull.Text)]"

# NOW split off paragraph 1
$trFull2 = $shp.TextFrame.TextRange
$trFull2.InsertBefore("Cut numbers, `r")
Write-Output "post-split: [$($shp.TextFrame.TextRange.Text)]"
